$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 is a text code that looks numeric ("001"); force Text format first so
# Excel keeps it as a string with the leading zero instead of coercing it to 1
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"

$ws.Range("N2").Value = "2019-12-31 00:00:00"

$ws.Range("O2").Value = 630859851.48
$ws.Range("P2").Value = 7302555.71
$ws.Range("Q2").Value = 418321378.92
$ws.Range("R2").Value = 519.1283094967
$ws.Range("S2").Value = 10217604.01
$ws.Range("T2").Value = -58.2336493785
$ws.Range("U2").Value = 152096395.06
$ws.Range("V2").Value = 79.2480044393
$ws.Range("W2").Value = 108806504.56
$ws.Range("X2").Value = 80423917.98999999
$ws.Range("Y2").Value = -8.053399364900001
$ws.Range("Z2").Value = 5736614.13
$ws.Range("AA2").Value = 1326.288159473
$ws.Range("AB2").Value = 522053346.92
$ws.Range("AC2").Value = 357.6488827862
$ws.Range("AD2").Value = 192.4091901108
$ws.Range("AE2").Value = 7.0164789561
$ws.Range("AF2").Value = 556.7801534562
$ws.Range("AG2").Value = 17.2473338262
